$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the "peptide_abundance_study_variable.N." labels in column A
#     with "sumIntensity_N" labels (rows 2-9), keep "run"/"treatment" header
#     and the UPS1/UPS2 values in column B unchanged. ---

$ws.Range("A1").Value = "run"
$ws.Range("B1").Value = "treatment"

$ws.Range("A2").Value = "sumIntensity_1"
$ws.Range("B2").Value = "UPS1"

$ws.Range("A3").Value = "sumIntensity_2"
$ws.Range("B3").Value = "UPS1"

$ws.Range("A4").Value = "sumIntensity_3"
$ws.Range("B4").Value = "UPS1"

$ws.Range("A5").Value = "sumIntensity_4"
$ws.Range("B5").Value = "UPS1"

$ws.Range("A6").Value = "sumIntensity_5"
$ws.Range("B6").Value = "UPS2"

$ws.Range("A7").Value = "sumIntensity_6"
$ws.Range("B7").Value = "UPS2"

$ws.Range("A8").Value = "sumIntensity_7"
$ws.Range("B8").Value = "UPS2"

$ws.Range("A9").Value = "sumIntensity_8"
$ws.Range("B9").Value = "UPS2"

# --- Widen column A so the longer "sumIntensity_N" labels are fully visible ---
$ws.Columns.Item(1).ColumnWidth = 32.85546875

# --- Update the active selection on the sheet ---
$ws.Range("G6").Select()
